$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Rename the header in C1: the config key changes from an index-based
# column reference to a letter-based one.
$ws.Range("C1").Value = "DATE_COLUMN_LETTER"

# Fix the swapped file-name / sheet-name values and set the new
# date-column-letter value.
$ws.Range("A2").Value = "weather_data"
$ws.Range("B2").Value = "Sheet1"
$ws.Range("C2").Value = "A"

# Move/record the active selection as it was when the author saved.
$ws.Range("C7").Select()
